# Fix minor typo/formatting in the "graphs-topoSort-connComponents" slide:
# superscript the "T" in "GT" (transpose graph notation) on the bullet
# "Reversed edges in GT stop it visiting nodes in SCCs yet to be found".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the paragraph containing the target sentence.
$targetPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -like "*Reversed edges in GT stop it visiting nodes in SCCs yet to be found*") {
        $targetPara = $para
        break
    }
}

# Find the "T" immediately following "Reversed edges in G" within the paragraph
# and raise it as a superscript, leaving the rest of the text untouched.
$offset = $targetPara.Text.IndexOf("Reversed edges in G") + ("Reversed edges in G").Length
$tChar = $tr.Characters($targetPara.Start + $offset, 1)
$tChar.Font.BaselineOffset = 0.3
